# Applies the LOT2005.docx content reshuffle described by the commit diff.
# The edit moves whole blocks of paragraph text between paragraphs while each
# paragraph keeps its own formatting (style / bold / italic) - so this script
# performs the move via Find and Replace on each run's literal text.
#
# Several of the text blocks being moved are each other's find/replace
# targets (e.g. slot 1's destination text is slot 5's source text), so a
# direct one-shot find/replace pass would clobber data. Instead we route
# every moved block through a neutral placeholder token first (Phase 1),
# then fill in the real destination text from the placeholders (Phase 2).

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $null = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

# ---- Phase 1: move each source runs text behind a unique placeholder ----
Replace-Text 'Desenvolver habilidades no campo da bioquímica através da execução e discussão de atividades práticas de laboratório.' '@@SLOT1@@'
Replace-Text 'Develop skills in the field of biochemistry through the execution and discussion of practical laboratory activities' '@@SLOT2@@'
Replace-Text '427823 - Adriane Maria Ferreira Milagres^l' '@@SLOT3@@'
Replace-Text '5082401 - André Moreni Lopes' '@@SLOT4@@'
Replace-Text '01Determinações analíticas envolvendo medidas de pH e condutividade. 02 Solubilização de bioativos hidrofóbicos. 03 soluções tampão. 04 caracterização de aminoácidos. 05 determinação de proteínas. 06  Separação e caracterização de proteinas. 07 cinética enzimática 08. dosagem de glicídios redutores 09 Caracterização de lipídeos.' '@@SLOT5@@'
Replace-Text '01 Analytical determinations involving pH and conductivity measurements. 02 Solubilization of hydrophobic bioactives. 03 Buffer solutions. 04 Characterization of amino acids. 05 Protein determination. 06 Separation and characterization of proteins. 07 Enzyme kinetics. 08 Assay of reducing sugars. 09 Characterization of lipids' '@@SLOT6@@'
Replace-Text '01 Uso de medidas de pH e condutividade. 02 Uso de tampões nas analises bioquimicas. 03 Demonstração da difusão seletiva de partículas de soluto através de membranas. 04 Solubilização de bioativos hidrofóbicos em ambiente aquoso através de tensoativos. 05  Aminoácidos: reação com ninhidrina. 06 Varredura de espectro: escolha do comprimento de onda ideal para dosagens fotométricas, relação absorbância e concentração; varredura de espectro de aminoácidos. 07 Cromatografia de aminoácidos: fundamentos físico, análise das fases envolvidas, escolha do solvente, solubilidade relativa dos componentes da amostra entre as duas fases, determinação de Rf;. 08 Fundamento químico das revelações de aminoácidos e proteínas. 9 Proteínas : separação e caracterização. 10 Dosagem e curva de proteína pelo método de Biureto; fundamento químico; obtenção da curva e aplicação prática da curva; 11. Enzimas: efeito do tempo.  curva de progresso, traçado e análise da curva, conceito de velocidade inicial, 12. influência da temperatura, aplicação da equação de Arrhenius, conceito de energia de ativação (traçado e análise da curva). 13. Dosagem de glicídios redutores: dosagem e curva padrão de glicídeos redutores pelo método do ácido 3,5-dinitro salicílico. 14 Caracterização de trigliceridos do oleo vegetal e determinação de alguns indices.' '@@SLOT7@@'
Replace-Text '1 Prova escrita (E), solicitando interpretação e análise de resultados experimentais. 2. cada grupo deverá compor um relatório em forma de seminário que deverá ser entregue para o professor e apresentado para a turma. A Avaliação será realizada por meio da participação nas aulas, relatório (R) e prova escrita (E).^l' '@@SLOT8@@'
Replace-Text 'A média final será calculada segundo a equação abaixo: ^lMédia final = 0,4 R + 0,6 P^l' '@@SLOT9@@'
Replace-Text 'A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculada pela equação: MF = (NF + PR)/2.' '@@SLOT10@@'
Replace-Text '1. Nelson, D.L., Cox, M.M. Princípios de bioquímica de Lehninger. Artmed^lEditora, 2022.^l2. Segel, I.H. Bioquímica Teoria e Problemas, São Paulo: Livros técnicos e Científicos Editora S.A, 1979.^l 3. Artigos e revisões da literatura ou outra bibliografia indicada no cronograma anual da disciplina.' '@@SLOT11@@'

# ---- Phase 2: drop in the real destination text for each placeholder ----
Replace-Text '@@SLOT1@@' '01Determinações analíticas envolvendo medidas de pH e condutividade. 02 Solubilização de bioativos hidrofóbicos. 03 soluções tampão. 04 caracterização de aminoácidos. 05 determinação de proteínas. 06  Separação e caracterização de proteinas. 07 cinética enzimática 08. dosagem de glicídios redutores 09 Caracterização de lipídeos.'
Replace-Text '@@SLOT2@@' '01 Analytical determinations involving pH and conductivity measurements. 02 Solubilization of hydrophobic bioactives. 03 Buffer solutions. 04 Characterization of amino acids. 05 Protein determination. 06 Separation and characterization of proteins. 07 Enzyme kinetics. 08 Assay of reducing sugars. 09 Characterization of lipids'
Replace-Text '@@SLOT3@@' 'Desenvolver habilidades no campo da bioquímica através da execução e discussão de atividades práticas de laboratório.^l'
Replace-Text '@@SLOT4@@' '01 Uso de medidas de pH e condutividade. 02 Uso de tampões nas analises bioquimicas. 03 Demonstração da difusão seletiva de partículas de soluto através de membranas. 04 Solubilização de bioativos hidrofóbicos em ambiente aquoso através de tensoativos. 05  Aminoácidos: reação com ninhidrina. 06 Varredura de espectro: escolha do comprimento de onda ideal para dosagens fotométricas, relação absorbância e concentração; varredura de espectro de aminoácidos. 07 Cromatografia de aminoácidos: fundamentos físico, análise das fases envolvidas, escolha do solvente, solubilidade relativa dos componentes da amostra entre as duas fases, determinação de Rf;. 08 Fundamento químico das revelações de aminoácidos e proteínas. 9 Proteínas : separação e caracterização. 10 Dosagem e curva de proteína pelo método de Biureto; fundamento químico; obtenção da curva e aplicação prática da curva; 11. Enzimas: efeito do tempo.  curva de progresso, traçado e análise da curva, conceito de velocidade inicial, 12. influência da temperatura, aplicação da equação de Arrhenius, conceito de energia de ativação (traçado e análise da curva). 13. Dosagem de glicídios redutores: dosagem e curva padrão de glicídeos redutores pelo método do ácido 3,5-dinitro salicílico. 14 Caracterização de trigliceridos do oleo vegetal e determinação de alguns indices.'
Replace-Text '@@SLOT5@@' '1 Prova escrita (E), solicitando interpretação e análise de resultados experimentais. 2. cada grupo deverá compor um relatório em forma de seminário que deverá ser entregue para o professor e apresentado para a turma. A Avaliação será realizada por meio da participação nas aulas, relatório (R) e prova escrita (E).'
Replace-Text '@@SLOT6@@' 'Develop skills in the field of biochemistry through the execution and discussion of practical laboratory activities'
Replace-Text '@@SLOT7@@' 'A média final será calculada segundo a equação abaixo: ^lMédia final = 0,4 R + 0,6 P'
Replace-Text '@@SLOT8@@' 'A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculada pela equação: MF = (NF + PR)/2.^l'
Replace-Text '@@SLOT9@@' '1. Nelson, D.L., Cox, M.M. Princípios de bioquímica de Lehninger. Artmed^lEditora, 2022.^l2. Segel, I.H. Bioquímica Teoria e Problemas, São Paulo: Livros técnicos e Científicos Editora S.A, 1979.^l 3. Artigos e revisões da literatura ou outra bibliografia indicada no cronograma anual da disciplina.^l'
Replace-Text '@@SLOT10@@' '427823 - Adriane Maria Ferreira Milagres'
Replace-Text '@@SLOT11@@' '5082401 - André Moreni Lopes'
